# Generate Report for Handoff
# Update the "Latest Handoff Datetime" for af22f8cb-2408-4abf-8703-dc848a9921b7.md
# on the "zh-cn" sheet, recording the newly generated handoff xliff timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("zh-cn")

# Row 5 = af22f8cb-2408-4abf-8703-dc848a9921b7.md, column H = "Latest Handoff Datetime"
$ws.Range("H5").Value = "2016-09-02 12:48:18"
